$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the harvestDate typo: "9.10.19" -> "09.10.19" (column A, rows 2-37) ---
# Build the corrected text via a formula first (so the engine treats it as
# literal text rather than auto-coercing the dotted string into a date),
# then convert the formulas to static values in one shot.
$dateRange = $ws.Range("A2:A37")
$dateRange.Formula = "=""09.10.19"""
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Normalize the harvester column (B) font back to the workbook's base font ---
$harvesterRange = $ws.Range("B2:B37")
$harvesterRange.Font.Name = "Calibri"
$harvesterRange.Font.Size = 12
$harvesterRange.Font.Color = -16777216

# --- Row heights: existing data rows 3-37 shrink to 15pt (row 2 keeps 16pt) ---
$ws.Range("A3:A37").EntireRow.RowHeight = 15

# --- Add the new trailing blank row 43 (height 15), extending the sheet's dimension ---
$ws.Rows.Item(43).RowHeight = 15
$ws.Cells.Item(43, 1).NumberFormat = "General"

# --- Update the active selection to G10 ---
$ws.Range("G10").Select()
